# Insert a new weekly data row for "Brócoli" (Macroferia Regional de Talca)
# above the existing row 432. Excel's row Insert shifts the former rows
# 432..479 down to 433..480, carrying all of their data (and formatting)
# along with them - exactly matching the target diff, which shows every
# row from 433 to 480 taking on the values that used to belong to the row
# immediately above it, while row 432 receives brand-new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("432:432").Insert()

$ws.Range("A432").Value = 5
$ws.Range("B432").Value = "Macroferia Regional de Talca"
$ws.Range("C432").Value = "Maule"
$ws.Range("D432").Value = 44946
$ws.Range("E432").Value = 7
$ws.Range("F432").Value = 100112023
$ws.Range("G432").Value = "Brócoli"
$ws.Range("H432").Value = "Sin especificar"
$ws.Range("I432").Value = "Primera"
$ws.Range("J432").Value = 2000
$ws.Range("K432").Value = 1000
$ws.Range("L432").Value = 1000
$ws.Range("M432").Value = 1000
$ws.Range("N432").Value = "`$/unidad"
$ws.Range("O432").Value = "Región del Maule"
$ws.Range("P432").Value = 1000
$ws.Range("Q432").Value = 1
$ws.Range("R432").Value = "Hortaliza"
